function Set-ParaText($Shape, $Index, $Value) {
    # Clear first, then re-fetch the paragraph and assign the final value.
    # Doing the replacement in one step can make the host preserve a
    # trailing character shared between old/new text as a leftover run;
    # clearing first avoids that and keeps a single <a:r> per paragraph.
    $range = $Shape.TextFrame.TextRange
    $range.Paragraphs($Index, 1).Text = ""
    $Shape.TextFrame.TextRange.Paragraphs($Index, 1).Text = $Value
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)

# Insert the new lead-in paragraph right after the existing (leading) blank
# paragraph, by appending it to that blank paragraph separated by a CR. This
# makes the new paragraph share the blank paragraph's "no explicit <a:pPr/>"
# formatting, matching the target rather than the pPr-bearing body bullets.
$tr = $shape.TextFrame.TextRange
$blankPara = $tr.Paragraphs(1, 1)
[void]$blankPara.InsertAfter([char]13 + "Compute Metrics and Draft Summary:")

# Rewrite the text of each existing metric/highlight paragraph in place so
# each one keeps its own <a:pPr/>. Paragraph indices shift down by one
# (the new paragraph became #2), so the metric bullets are now #3-#8.
Set-ParaText $shape 3 "**Total GLA**: 313,219 m² (222,221 m² for Ingram Micro + 90,998 m² for CNH Industrial)"
Set-ParaText $shape 4 "**Occupancy**: Potentially 0% (leases for both tenants have expired)"
Set-ParaText $shape 5 "**WALT**: 0 years (Weighted Average Lease Term is 0 due to expired leases)"
Set-ParaText $shape 6 "**In-Place Rent**: Not applicable (leases have expired)"
Set-ParaText $shape 7 "**Key Highlight 1**: Strategic logistics location in Daventry, UK, with proximity to major transportation routes."
Set-ParaText $shape 8 "**Key Highlight 2**: Significant leasable area with potential for stable cash flow if leases are renewed or new tenants are secured."
